$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Turn off the AutoFilter - removes the <autoFilter> element from the sheet
$ws.AutoFilterMode = $false

# 2. Remove the now-orphaned "_xlnm._FilterDatabase" defined name so
#    <definedNames> ends up empty
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# 3. Fix the wrong SIRET number
$ws.Range("R2").Value = "21920044100018"

# 4. Slightly narrow the data columns (A:AD)
$ws.Columns("A:AD").ColumnWidth = 7.8

# 5. Give the header row and the data row an explicit row height
$ws.Rows("1:2").RowHeight = 14.25

# 6. C1 and R1 used the "text" number-format style; switch them back to the
#    plain style already used by A1/B1, without touching their text value
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A1").Copy()
$ws.Range("R1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
